$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that switch from STAY -> CHARTER
$toCharter = @(63, 64, 83, 84, 85, 88, 89, 91)
foreach ($r in $toCharter) {
    $ws.Cells.Item($r, 1).Value = "CHARTER"
}

# Rows that switch from CHARTER -> STAY
$toStay = @(110, 114, 115, 116)
foreach ($r in $toStay) {
    $ws.Cells.Item($r, 1).Value = "STAY"
}
